$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of the existing header cell (H1) onto the new header
# cells so they match the look (bold, bordered, centered) of the rest
# of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats

# Header cells for the two new columns, I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for column I (I0) and column J (IF), rows 2 through 51
$iVals = @(7,9,8,7,7,6,7,7,5,6,8,7,7,8,7,6,9,5,7,7,7,7,6,9,5,9,8,7,5,9,8,7,7,8,6,6,9,6,7,5,11,8,7,6,8,7,8,7,5,4)
$jVals = @(7,9,8,8,7,6,7,8,6,7,8,7,8,8,7,6,9,5,7,7,8,7,7,10,6,9,8,8,5,9,9,8,8,9,6,7,9,7,7,6,11,8,7,7,8,7,8,7,5,4)

for ($i = 0; $i -lt $iVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$i]
    $ws.Cells.Item($row, 10).Value = $jVals[$i]
}
